# Refresh the crypto price/volume data (and two reordered rows) to match the
# latest scrape, as produced by the "Updated cryptos list" GitHub Action.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row updates: RowNumber -> @{ Col = NewValue }  (columns: B=Coin, C=Link, D=Price, E=Volume(1h))
$updates = [ordered]@{
    2 = @{ D='27.690.21'; E='  -0.66%  ' }
    3 = @{ D='1.894.17'; E='  +1.21%  ' }
    4 = @{ D='1.001'; E='  -1.14%  ' }
    5 = @{ D='312.72'; E='  -0.34%  ' }
    6 = @{ E='  -1.07%  ' }
    7 = @{ D='0.4936'; E='  +1.92%  ' }
    8 = @{ D='0.3803'; E='  -0.63%  ' }
    9 = @{ D='0.07327'; E='  -0.59%  ' }
    10 = @{ D='0.9140'; E='  -2.79%  ' }
    11 = @{ D='20.56'; E='  -2.25%  ' }
    12 = @{ B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='1.931.63'; E='  +3.04%  ' }
    13 = @{ B='TRON'; C='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D='0.07674'; E='  -1.90%  ' }
    14 = @{ D='5.474'; E='  -0.47%  ' }
    15 = @{ D='6.593'; E='  -0.36%  ' }
    16 = @{ D='91.07'; E='  -0.17%  ' }
    17 = @{ D='1.002'; E='  -1.10%  ' }
    18 = @{ D='0.000008775'; E='  -1.16%  ' }
    20 = @{ D='27.763.45'; E='  -0.47%  ' }
    21 = @{ D='14.52'; E='  -2.28%  ' }
    22 = @{ D='5.125'; E='  -0.03%  ' }
    23 = @{ D='2.146.11'; E='  +1.43%  ' }
    24 = @{ E='  -1.11%  ' }
    25 = @{ D='1.909'; E='  -2.04%  ' }
    26 = @{ D='153.44'; E='  -2.29%  ' }
    27 = @{ D='18.38'; E='  -1.07%  ' }
    28 = @{ D='2.143'; E='  +4.18%  ' }
    29 = @{ D='115.61'; E='  -0.37%  ' }
    30 = @{ D='4.894'; E='  -1.92%  ' }
    31 = @{ D='0.08936'; E='  +0.18%  ' }
    32 = @{ D='3.198'; E='  -4.04%  ' }
    33 = @{ E='  -0.90%  ' }
    34 = @{ D='0.7655'; E='  -0.60%  ' }
    35 = @{ E='  -0.43%  ' }
    36 = @{ D='0.02030'; E='  -0.98%  ' }
    37 = @{ D='2.523'; E='  -7.83%  ' }
    38 = @{ E='  -3.56%  ' }
    39 = @{ D='0.05279'; E='  -1.65%  ' }
    40 = @{ D='0.5469'; E='  -2.55%  ' }
    41 = @{ E='  -0.56%  ' }
    42 = @{ D='6.908'; E='  -2.04%  ' }
    43 = @{ D='8.526'; E='  -0.65%  ' }
    44 = @{ B='Quant'; C='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D='112.46'; E='  +6.61%  ' }
    45 = @{ B='Algorand'; C='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D='0.1516'; E='  -1.23%  ' }
    46 = @{ E='  -1.46%  ' }
    47 = @{ D='0.4791'; E='  -1.91%  ' }
    48 = @{ E='  -1.21%  ' }
    49 = @{ D='1.630'; E='  -2.36%  ' }
    50 = @{ D='67.48'; E='  -0.99%  ' }
    51 = @{ D='0.06049' }
}

$colIndex = @{ B = 2; C = 3; D = 4; E = 5 }

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    # Column D holds numeric-looking text (e.g. "0.9140", "1.001") that must stay
    # text so formatting/trailing zeros and thousands-dot groupings are preserved.
    if ($rowData.Contains("D")) {
        $ws.Cells.Item($row, $colIndex.D).NumberFormat = "@"
    }
    foreach ($col in $rowData.Keys) {
        $ws.Cells.Item($row, $colIndex[$col]).Value = $rowData[$col]
    }
    if ($rowData.Contains("D")) {
        $ws.Cells.Item($row, $colIndex.D).ClearFormats()
    }
}

Write-Host "Applied cryptos update"